$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet grows from 10 columns (A:J) to 13 columns (A:M). Insert 3 blank
# columns at the tail (after J) so we have room for the new L/M columns, then
# rewrite every cell (headers + data) with its final value, since virtually
# every column's content (not just its position) changes.
# ---------------------------------------------------------------------------
$ws.Columns.Item(11).Insert()
$ws.Columns.Item(11).Insert()
$ws.Columns.Item(11).Insert()

# ---------------------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------------------
$headers = @("Distance","Source_Type","Outcome_Var","Exclude_Touching","Include_area","Include_population","Coefficient","CI_lower","CI_upper","P-value","F-statistic","n_observations","R-squared")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = $headers[$i]
}

# Re-apply the bold/centered/bordered header style (already on A1) to the
# newly added header cells K1:M1 so every header cell shares style index 1.
$ws.Cells.Item(1,1).Copy()
$ws.Range("K1:M1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Data rows 2-5
# ---------------------------------------------------------------------------
# Columns: A Distance | B Source_Type | C Outcome_Var | D Exclude_Touching |
#          E Include_area | F Include_population | G Coefficient | H CI_lower |
#          I CI_upper | J P-value | K F-statistic | L n_observations | M R-squared
$data = @(
    @{A=100; B=""; C="num_power_stations"; D=$true; E=$true; F=$true; G=1.041202766831799;  H=-0.01699108037007502; I=2.099396614033673; J=0.05378687140079773;  K=50.27570070516628; L=1065; M=0.12446255893452},
    @{A=150; B=""; C="num_power_stations"; D=$true; E=$true; F=$true; G=1.345818573717814;  H=0.3504478116466351;   I=2.341189335788992; J=0.008082660657925679; K=48.5934357513277;  L=1479; M=0.08994452030390987},
    @{A=200; B=""; C="num_power_stations"; D=$true; E=$true; F=$true; G=1.4267694629089;    H=0.4495863445585357;   I=2.403952581259264; J=0.004238403492816602; K=39.61248866532684; L=1672; M=0.06650714922733092},
    @{A=250; B=""; C="num_power_stations"; D=$true; E=$true; F=$true; G=1.533055607148542;  H=0.6003052936548631;   I=2.465805920642221; J=0.001288495079580839; K=40.27802394708902; L=1853; M=0.06134225900982593}
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M")
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowData = $data[$r]
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$rowNum")
        $value = $rowData[$col]
        if ($col -eq "B") {
            # Source_Type is blank for every row, but it must remain a *text*
            # cell (matching the sheet's empty-string convention used
            # elsewhere, e.g. Outcome_Var on the original sheet) rather than
            # an untyped/numeric blank. A leading quote forces Excel to store
            # it as text; resetting the style afterwards drops the stray
            # "quote prefix" formatting that introduces.
            $cell.Value2 = "'"
            $cell.Style = "Normal"
        } else {
            $cell.Value2 = $value
        }
    }
}
